$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.813.63"
$ws.Range("D2").Style = $cellStyle
$ws.Range("E2").Value = "  +2.59%  "

$cellStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.092.16"
$ws.Range("D3").Style = $cellStyle
$ws.Range("E3").Value = "  +2.37%  "

$ws.Range("E4").Value = "  +0.06%  "

$cellStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.35"
$ws.Range("D5").Style = $cellStyle
$ws.Range("E5").Value = "  +0.37%  "

$ws.Range("E6").Value = "  +0.86%  "

$cellStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.70"
$ws.Range("D7").Style = $cellStyle
$ws.Range("E7").Value = "  +1.61%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  +1.82%  "

$ws.Range("E10").Value = "  +0.11%  "

$cellStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.104"
$ws.Range("D11").Style = $cellStyle
$ws.Range("E11").Value = "  -0.56%  "

$cellStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.403.01"
$ws.Range("D12").Style = $cellStyle
$ws.Range("E12").Value = "  +2.59%  "

$cellStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.99"
$ws.Range("D13").Style = $cellStyle
$ws.Range("E13").Value = "  +3.56%  "

$ws.Range("E14").Value = "  +4.58%  "

$ws.Range("E15").Value = "  +3.20%  "

$cellStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.47"
$ws.Range("D16").Style = $cellStyle
$ws.Range("E16").Value = "  -0.69%  "

$cellStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.096.43"
$ws.Range("D17").Style = $cellStyle
$ws.Range("E17").Value = "  +2.69%  "

$cellStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "38.739.45"
$ws.Range("D18").Style = $cellStyle
$ws.Range("E18").Value = "  +2.61%  "

$cellStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.86"
$ws.Range("D19").Style = $cellStyle
$ws.Range("E19").Value = "  +3.34%  "

$cellStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.02"
$ws.Range("D20").Style = $cellStyle
$ws.Range("E20").Value = "  +1.98%  "

$ws.Range("E21").Value = "  +1.43%  "

$cellStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.44"
$ws.Range("D22").Style = $cellStyle
$ws.Range("E22").Value = "  +1.13%  "

$ws.Range("E23").Value = "  -0.38%  "

$cellStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.42"
$ws.Range("D24").Style = $cellStyle
$ws.Range("E24").Value = "  -0.42%  "

$cellStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.34"
$ws.Range("D25").Style = $cellStyle
$ws.Range("E25").Value = "  +2.41%  "

$cellStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.64"
$ws.Range("D26").Style = $cellStyle
$ws.Range("E26").Value = "  +0.98%  "

$cellStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.39"
$ws.Range("D27").Style = $cellStyle
$ws.Range("E27").Value = "  +0.15%  "

$cellStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.136"
$ws.Range("D28").Style = $cellStyle
$ws.Range("E28").Value = "  +6.55%  "

$cellStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.45"
$ws.Range("D29").Style = $cellStyle
$ws.Range("E29").Value = "  +12.95%  "

$cellStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.16"
$ws.Range("D30").Style = $cellStyle
$ws.Range("E30").Value = "  +1.89%  "

$ws.Range("E31").Value = "  +0.46%  "

$ws.Range("E32").Value = "  +4.27%  "

$ws.Range("E33").Value = "  +1.89%  "

$ws.Range("E34").Value = "  +4.62%  "

$ws.Range("E35").Value = "  +1.25%  "

$ws.Range("B36").Value = "THORChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$cellStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.44"
$ws.Range("D36").Style = $cellStyle
$ws.Range("E36").Value = "  -1.35%  "

$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$cellStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.39"
$ws.Range("D37").Style = $cellStyle
$ws.Range("E37").Value = "  +2.07%  "

$ws.Range("E38").Value = "  +2.35%  "

$ws.Range("E39").Value = "  +0.12%  "

$cellStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.25"
$ws.Range("D40").Style = $cellStyle
$ws.Range("E40").Value = "  +1.07%  "

$cellStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "101.19"
$ws.Range("D41").Style = $cellStyle
$ws.Range("E41").Value = "  +3.82%  "

$cellStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.538.36"
$ws.Range("D42").Style = $cellStyle
$ws.Range("E42").Value = "  +0.65%  "

$cellStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0223"
$ws.Range("D43").Style = $cellStyle
$ws.Range("E43").Value = "  +3.39%  "

$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$cellStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0925"
$ws.Range("D44").Style = $cellStyle
$ws.Range("E44").Value = "  +2.24%  "

$ws.Range("B45").Value = "HuobiToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$cellStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.82"
$ws.Range("D45").Style = $cellStyle
$ws.Range("E45").Value = "  -0.74%  "

$ws.Range("E46").Value = "  +8.40%  "

$cellStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.09"
$ws.Range("D48").Style = $cellStyle
$ws.Range("E48").Value = "  -2.96%  "

$ws.Range("E49").Value = "  +2.38%  "

$ws.Range("E50").Value = "  +0.86%  "

$cellStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.288.99"
$ws.Range("D51").Style = $cellStyle
$ws.Range("E51").Value = "  +2.49%  "
